$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(13, 8).Value = 18472.5
$ws.Cells.Item(13, 10).Value = 7963.3335
$ws.Cells.Item(13, 12).Value = 7963.3335
$ws.Cells.Item(13, 14).Value = -8301.333500000001

$ws.Cells.Item(62, 8).Value = 5729.8
$ws.Cells.Item(62, 9).Value = 4333
$ws.Cells.Item(62, 10).Value = 7825
$ws.Cells.Item(62, 11).Value = 4333
$ws.Cells.Item(62, 12).Value = 7825
$ws.Cells.Item(62, 13).Value = -3709
$ws.Cells.Item(62, 14).Value = -9073

$ws.Cells.Item(65, 8).Value = 5729.8
$ws.Cells.Item(65, 9).Value = 4333
$ws.Cells.Item(65, 10).Value = 7825
$ws.Cells.Item(65, 11).Value = 21665
$ws.Cells.Item(65, 12).Value = 39125
$ws.Cells.Item(65, 13).Value = -18545
$ws.Cells.Item(65, 14).Value = -45365

$ws.Cells.Item(100, 8).Value = 2239.5386
$ws.Cells.Item(100, 9).Value = 2118.25
$ws.Cells.Item(100, 10).Value = 3695
$ws.Cells.Item(100, 11).Value = 2118.25
$ws.Cells.Item(100, 12).Value = 3695
$ws.Cells.Item(100, 13).Value = -1577.25
$ws.Cells.Item(100, 14).Value = -4777

$ws.Cells.Item(111, 8).Value = 9579.5
$ws.Cells.Item(111, 10).Value = 27800
$ws.Cells.Item(111, 12).Value = 83400
$ws.Cells.Item(111, 14).Value = -89534

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 601.8823
$ws.Cells.Item(2, 9).Value = 296.0909
$ws.Cells.Item(2, 11).Value = 296.0909
$ws.Cells.Item(2, 13).Value = -183.0909

$ws.Cells.Item(32, 8).Value = 7611.8887
$ws.Cells.Item(32, 9).Value = 3162.45
$ws.Cells.Item(32, 10).Value = 26346.37
$ws.Cells.Item(32, 11).Value = 3162.45
$ws.Cells.Item(32, 12).Value = 26346.37
$ws.Cells.Item(32, 13).Value = -2875.45
$ws.Cells.Item(32, 14).Value = -26920.37

$ws.Cells.Item(61, 8).Value = 3906.1155
$ws.Cells.Item(61, 10).Value = 5562.4585
$ws.Cells.Item(61, 12).Value = 5562.4585
$ws.Cells.Item(61, 14).Value = -5986.4585

$ws.Cells.Item(74, 8).Value = 35181.97
$ws.Cells.Item(74, 9).Value = 41034.64
$ws.Cells.Item(74, 10).Value = 16892.375
$ws.Cells.Item(74, 11).Value = 41034.64
$ws.Cells.Item(74, 12).Value = 16892.375
$ws.Cells.Item(74, 13).Value = -40160.64
$ws.Cells.Item(74, 14).Value = -18640.375

$ws.Cells.Item(77, 8).Value = 35181.97
$ws.Cells.Item(77, 9).Value = 41034.64
$ws.Cells.Item(77, 10).Value = 16892.375
$ws.Cells.Item(77, 11).Value = 205173.2
$ws.Cells.Item(77, 12).Value = 84461.875
$ws.Cells.Item(77, 13).Value = -200805.2
$ws.Cells.Item(77, 14).Value = -93197.875

$ws.Cells.Item(104, 8).Value = 74706.336
$ws.Cells.Item(104, 10).Value = 74706.336
$ws.Cells.Item(104, 12).Value = 74706.336
$ws.Cells.Item(104, 14).Value = -81694.336

$ws.Cells.Item(116, 8).Value = 601.8823
$ws.Cells.Item(116, 9).Value = 296.0909
$ws.Cells.Item(116, 11).Value = 296.0909
$ws.Cells.Item(116, 13).Value = 1997.9091

$ws.Cells.Item(136, 8).Value = 3906.1155
$ws.Cells.Item(136, 10).Value = 5562.4585
$ws.Cells.Item(136, 12).Value = 16687.3755
$ws.Cells.Item(136, 14).Value = -21787.3755

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 601.8823
$ws.Cells.Item(3, 9).Value = 296.0909
$ws.Cells.Item(3, 11).Value = 296.0909
$ws.Cells.Item(3, 13).Value = -182.0909

$ws.Cells.Item(99, 8).Value = 122175.414
$ws.Cells.Item(99, 9).Value = 85474.336
$ws.Cells.Item(99, 11).Value = 85474.336
$ws.Cells.Item(99, 13).Value = -83976.336

$ws.Cells.Item(103, 8).Value = 13999.8
$ws.Cells.Item(103, 10).Value = 13999.8
$ws.Cells.Item(103, 12).Value = 13999.8
$ws.Cells.Item(103, 14).Value = -16343.8

$ws.Cells.Item(105, 8).Value = 2164.6924
$ws.Cells.Item(105, 9).Value = 2184.4666
$ws.Cells.Item(105, 10).Value = 2098.7778
$ws.Cells.Item(105, 11).Value = 2184.4666
$ws.Cells.Item(105, 12).Value = 2098.7778
$ws.Cells.Item(105, 13).Value = -437.4666000000002
$ws.Cells.Item(105, 14).Value = -5592.7778

$ws.Cells.Item(106, 8).Value = 30000
$ws.Cells.Item(106, 10).Value = 30000
$ws.Cells.Item(106, 12).Value = 30000
$ws.Cells.Item(106, 14).Value = -32524

$ws.Cells.Item(134, 8).Value = 2658.9583
$ws.Cells.Item(134, 9).Value = 2509.7727
$ws.Cells.Item(134, 10).Value = 4300
$ws.Cells.Item(134, 11).Value = 7529.3181
$ws.Cells.Item(134, 12).Value = 12900
$ws.Cells.Item(134, 13).Value = -4994.3181
$ws.Cells.Item(134, 14).Value = -17970

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(88, 8).Value = 21138.334
$ws.Cells.Item(88, 10).Value = 15171
$ws.Cells.Item(88, 12).Value = 15171
$ws.Cells.Item(88, 14).Value = -15983

$ws.Cells.Item(91, 8).Value = 21138.334
$ws.Cells.Item(91, 10).Value = 15171
$ws.Cells.Item(91, 12).Value = 15171
$ws.Cells.Item(91, 14).Value = -17979

$ws.Cells.Item(105, 8).Value = 1160.8823
$ws.Cells.Item(105, 9).Value = 1299.091
$ws.Cells.Item(105, 10).Value = 907.5
$ws.Cells.Item(105, 11).Value = 1299.091
$ws.Cells.Item(105, 12).Value = 907.5
$ws.Cells.Item(105, 13).Value = 447.9090000000001
$ws.Cells.Item(105, 14).Value = -4401.5

$ws.Cells.Item(122, 8).Value = 2186.35
$ws.Cells.Item(122, 10).Value = 3524.625
$ws.Cells.Item(122, 12).Value = 10573.875
$ws.Cells.Item(122, 14).Value = -15473.875

$ws.Cells.Item(141, 8).Value = 396568
$ws.Cells.Item(141, 10).Value = 429275.16
$ws.Cells.Item(141, 12).Value = 429275.16
$ws.Cells.Item(141, 14).Value = -439635.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 564.38464
$ws.Cells.Item(5, 9).Value = 564.38464
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 1693.15392
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = -1581.15392
$ws.Cells.Item(5, 14).ClearContents()

$ws.Cells.Item(37, 8).Value = 114998.125
$ws.Cells.Item(37, 10).Value = 114998.125
$ws.Cells.Item(37, 12).Value = 344994.375
$ws.Cells.Item(37, 14).Value = -345218.375

$ws.Cells.Item(80, 8).Value = 200
$ws.Cells.Item(80, 10).Value = 200
$ws.Cells.Item(80, 12).Value = 600
$ws.Cells.Item(80, 14).Value = -2472

$ws.Cells.Item(83, 8).Value = 200
$ws.Cells.Item(83, 10).Value = 200
$ws.Cells.Item(83, 12).Value = 1800
$ws.Cells.Item(83, 14).Value = -11160

$ws.Cells.Item(92, 8).Value = 1440.5
$ws.Cells.Item(92, 10).Value = 1440.5
$ws.Cells.Item(92, 12).Value = 4321.5
$ws.Cells.Item(92, 14).Value = -6817.5

$ws.Cells.Item(118, 8).Value = 3249.8333
$ws.Cells.Item(118, 9).Value = 3249.8333
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 11).Value = 9749.499899999999
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 13).Value = -8506.499899999999
$ws.Cells.Item(118, 14).ClearContents()

$ws.Cells.Item(135, 8).Value = 564.38464
$ws.Cells.Item(135, 9).Value = 564.38464
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 5079.46176
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 13).Value = -2544.46176
$ws.Cells.Item(135, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 100000980
$ws.Cells.Item(102, 9).Value = 1094
$ws.Cells.Item(102, 11).Value = 1094
$ws.Cells.Item(102, 13).Value = 528

$ws.Cells.Item(122, 8).Value = 2702.2974
$ws.Cells.Item(122, 9).Value = 2189.0386
$ws.Cells.Item(122, 10).Value = 3915.4546
$ws.Cells.Item(122, 11).Value = 6567.1158
$ws.Cells.Item(122, 12).Value = 11746.3638
$ws.Cells.Item(122, 13).Value = -4117.1158
$ws.Cells.Item(122, 14).Value = -16646.3638

$ws.Cells.Item(132, 8).Value = 3009.4517
$ws.Cells.Item(132, 9).Value = 2154.5652
$ws.Cells.Item(132, 10).Value = 5467.25
$ws.Cells.Item(132, 11).Value = 6463.6956
$ws.Cells.Item(132, 12).Value = 16401.75
$ws.Cells.Item(132, 13).Value = -3933.6956
$ws.Cells.Item(132, 14).Value = -21461.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1714.3103
$ws.Cells.Item(22, 9).Value = 1138.9231
$ws.Cells.Item(22, 10).Value = 2181.8125
$ws.Cells.Item(22, 11).Value = 1138.9231
$ws.Cells.Item(22, 12).Value = 2181.8125
$ws.Cells.Item(22, 13).Value = -843.9231
$ws.Cells.Item(22, 14).Value = -2771.8125

$ws.Cells.Item(27, 8).Value = 1714.3103
$ws.Cells.Item(27, 9).Value = 1138.9231
$ws.Cells.Item(27, 10).Value = 2181.8125
$ws.Cells.Item(27, 11).Value = 1138.9231
$ws.Cells.Item(27, 12).Value = 2181.8125
$ws.Cells.Item(27, 13).Value = -1031.9231
$ws.Cells.Item(27, 14).Value = -2395.8125

$ws.Cells.Item(122, 8).Value = 438663.34
$ws.Cells.Item(122, 9).Value = 628203.5600000001
$ws.Cells.Item(122, 10).Value = 5428.5713
$ws.Cells.Item(122, 11).Value = 1884610.68
$ws.Cells.Item(122, 12).Value = 16285.7139
$ws.Cells.Item(122, 13).Value = -1882160.68
$ws.Cells.Item(122, 14).Value = -21185.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(30, 8).Value = 14999.5
$ws.Cells.Item(30, 10).Value = 14999.5
$ws.Cells.Item(30, 12).Value = 14999.5
$ws.Cells.Item(30, 14).Value = -15213.5

$ws.Cells.Item(107, 8).Value = 26114.82
$ws.Cells.Item(107, 9).Value = 473.0345
$ws.Cells.Item(107, 11).Value = 1419.1035
$ws.Cells.Item(107, 13).Value = 500.8965000000001
